$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9-13 had their observation data rotated between rows (same
# formatting/columns, values shuffled). Each block below writes the new
# value for every cell that actually changes, per row.
#
# Columns I (Antal) and Y/AA (Start/Slutdatum) hold text that looks like a
# number/date ("1", "2023-08-11", ...) but must stay plain text, matching
# how the source file stores them. Force text format first (only on the
# cells whose value actually changes) so Excel does not auto-coerce the
# value into a real number/date.
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I10").NumberFormat = "@"
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA11").NumberFormat = "@"
$ws.Range("AA12").NumberFormat = "@"

# Row 9 (now holds what used to be row 12's record)
$ws.Range("A9").Value = 111611165
$ws.Range("B9").Value = 84741
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 37
$ws.Range("F9").Value = "Jättekamskivling"
$ws.Range("G9").Value = "Amanita ceciliae"
$ws.Range("H9").Value = "(Berk. & Broome) Bas"
$ws.Range("I9").Value = "1"
$ws.Range("Q9").Value = 663088.0668624006
$ws.Range("R9").Value = 6634684.960451891
$ws.Range("Y9").Value = "2023-08-11"
$ws.Range("AA9").Value = "2023-08-11"
$ws.Range("AC9").Value = "1 ex. under ek och hassel."
$ws.Range("AX9").Value = "Gillis Aronsson"

# Row 10 (now holds what used to be row 9's record)
$ws.Range("A10").Value = 111611138
$ws.Range("B10").Value = 81796
$ws.Range("E10").Value = 5406
$ws.Range("F10").Value = "Gulmjölkig storskål"
$ws.Range("G10").Value = "Peziza succosa"
$ws.Range("H10").Value = "Berk."
$ws.Range("I10").Value = "3"
$ws.Range("Q10").Value = 663213.3366271106
$ws.Range("R10").Value = 6634830.464506784
$ws.Range("Y10").Value = "2023-08-12"
$ws.Range("AA10").Value = "2023-08-12"
$ws.Range("AC10").Value = "3 ex. på bar jord och i lövförna."
$ws.Range("AX10").Value = "Gillis Aronsson, Cajsa Björkén"

# Row 11 (now holds what used to be row 10's record)
$ws.Range("A11").Value = 111611146
$ws.Range("B11").Value = 88630
$ws.Range("E11").Value = 4823
$ws.Range("F11").Value = "Hasselsopp"
$ws.Range("G11").Value = "Leccinellum pseudoscabrum"
$ws.Range("H11").Value = "(Kallenb.) Mikšík"
$ws.Range("Q11").Value = 663088.0668624006
$ws.Range("R11").Value = 6634684.960451891
$ws.Range("Y11").Value = "2023-08-11"
$ws.Range("AA11").Value = "2023-08-11"
$ws.Range("AC11").Value = "1 ex. under ek och hassel."
$ws.Range("AX11").Value = "Gillis Aronsson"

# Row 12 (now holds what used to be row 13's record)
$ws.Range("A12").Value = 111611145
$ws.Range("B12").Value = 88630
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 4823
$ws.Range("F12").Value = "Hasselsopp"
$ws.Range("G12").Value = "Leccinellum pseudoscabrum"
$ws.Range("H12").Value = "(Kallenb.) Mikšík"
$ws.Range("Q12").Value = 663143.8264147732
$ws.Range("R12").Value = 6634793.669287071
$ws.Range("Y12").Value = "2023-08-12"
$ws.Range("AA12").Value = "2023-08-12"
$ws.Range("AC12").Value = "1 ex. i lövförna under hassel."
$ws.Range("AX12").Value = "Gillis Aronsson, Cajsa Björkén"

# Row 13 (now holds what used to be row 11's record)
$ws.Range("A13").Value = 111611158
$ws.Range("B13").Value = 86021
$ws.Range("E13").Value = 4037
$ws.Range("F13").Value = "Bolmörtsskivling"
$ws.Range("G13").Value = "Entoloma sinuatum"
$ws.Range("H13").Value = "(Bull.) P.Kumm."
$ws.Range("Q13").Value = 663128.0992466732
$ws.Range("R13").Value = 6634761.25188593
$ws.Range("AC13").Value = "1 ex. i lövförna under ek och hassel."
